$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Cell value edits -----------------------------------------------------
# Non-statistical indicator / Use geocodes for mapping: "false" -> "False"
$ws.Range("B17").Value = "'False"
$ws.Range("B20").Value = "'False"

# Show source 1 / Show source 2: blank -> "True"
$ws.Range("B23").Value = "'True"
$ws.Range("B34").Value = "'True"

# --- Data validation --------------------------------------------------------
# Remove B17 and B20 from the old lowercase 'true/'false validation list
# (only B21 keeps the original rule).
$ws.Range("B17").Validation.Delete()
$ws.Range("B20").Validation.Delete()

# Add the new, capitalised 'True/'False validation rule to B17, B20, B23, B34
$ws.Range("B17").Validation.Add(3, 1, 1, """'True, 'False""")
$ws.Range("B20").Validation.Add(3, 1, 1, """'True, 'False""")
$ws.Range("B23").Validation.Add(3, 1, 1, """'True, 'False""")
$ws.Range("B34").Validation.Add(3, 1, 1, """'True, 'False""")

# --- Selection state --------------------------------------------------------
# Active cell moves from B21 to B34, and the view scrolls back to the top
# (no more frozen/top-left offset at A7).
[void]$ws.Range("B34").Select()
